# Report regenerated for a new source file GUID:
#   a7a731e2-3479-4062-b105-68d3de5b591e -> 2a18b13d-1bf0-43b7-9d65-720c7171a8eb
# (and the associated xliff hash a4bd988afe2d02ee69cb11e38aea39115e380a48 ->
#  d33127294674c68451d60c2cd45ee0eadf99f8a5), plus refreshed handoff/handback
# timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "a7a731e2-3479-4062-b105-68d3de5b591e"
$newGuid = "2a18b13d-1bf0-43b7-9d65-720c7171a8eb"
$oldHash = "a4bd988afe2d02ee69cb11e38aea39115e380a48"
$newHash = "d33127294674c68451d60c2cd45ee0eadf99f8a5"

# The external hyperlink target (unchanged by this edit - only the displayed
# text is refreshed to the new file name).
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8f7bf646e0211c3cc2b070276a988f8ae3fb17d9/e2e/$oldGuid.md"

# Original custom "Hyperlink" font formatting used throughout this workbook
# (single underline, cornflower blue FF6495ED == RGB(101,149,237)).
$linkColor = 15570276

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, "", "", "e2e\$newGuid.md")
$wsOverview.Range("B2").Font.Underline = $true
$wsOverview.Range("B2").Font.Color = $linkColor
$wsOverview.Range("G2").Value = "2016-10-20 00:40:20"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $linkAddress, "", "", "$newGuid.md")
$wsZh.Range("A2").Font.Underline = $true
$wsZh.Range("A2").Font.Color = $linkColor
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-10-20 00:40:08"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $linkAddress, "", "", "$newGuid.md")
$wsDe.Range("A2").Font.Underline = $true
$wsDe.Range("A2").Font.Color = $linkColor
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-10-20 00:40:20"
